# Add files via upload
# The "Survey 4" block (rows 16-22, column C) was missing household counts.
# This fills in the counts that accompany the existing Survey/Source rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C16").Value = 31
$ws.Range("C17").Value = 27
$ws.Range("C18").Value = 3
$ws.Range("C19").Value = 20
$ws.Range("C20").Value = 65
$ws.Range("C21").Value = 42
$ws.Range("C22").Value = 0

# Reflect the author's final view/selection state on the sheet.
$null = $ws.Range("A15").Select()
$null = $ws.Range("F26").Select()
